$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: fill in contribution data for the 2nd contributor (mirrors row 2's pattern) ---

# B3: date 2024-02-02, reuse B2's existing format (style s="6", numFmtId 14) instead of
# minting a brand-new number format, so we pick up the style via a format-only copy/paste.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("B3").Value = "2/2/2024"

# C3:F3 and H3: 16.7% contribution each - these cells already carry the workbook's
# 0.0% percentage style (s="2"), so a plain value write keeps that style intact.
$ws.Range("C3").Value = 0.167
$ws.Range("D3").Value = 0.167
$ws.Range("E3").Value = 0.167
$ws.Range("F3").Value = 0.167
$ws.Range("H3").Value = 0.167

# G3: same 16.7% contribution, but tagged with the built-in "Percent" cell style
# (new cellStyleXf / cellStyle entries) rather than the plain 0.0% style.
$ws.Range("G3").Value = 0.167
$ws.Range("G3").Style = "Percent"

# J3: total of 100% for the row, entered as a plain value (not a SUM formula).
$ws.Range("J3").Value = 1

# --- Selection moves to C4 ---
$ws.Range("C4").Select() | Out-Null
